# Update F-column (participant/view-count) figures across the four sheets
# of the 杭州-漫展信息 workbook, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 12950
$ws.Range("F3").Value  = 7263
$ws.Range("F7").Value  = 1025
$ws.Range("F8").Value  = 160
$ws.Range("F9").Value  = 373
$ws.Range("F10").Value = 1055
$ws.Range("F11").Value = 18
$ws.Range("F15").Value = 277
$ws.Range("F18").Value = 292
$ws.Range("F19").Value = 324
$ws.Range("F21").Value = 264
$ws.Range("F23").Value = 5322
$ws.Range("F26").Value = 331
$ws.Range("F27").Value = 2036
$ws.Range("F28").Value = 105
$ws.Range("F30").Value = 1404
$ws.Range("F35").Value = 3758

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value  = 84
$ws.Range("F14").Value = 82
$ws.Range("F18").Value = 49

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2075

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 2075
$ws.Range("F5").Value  = 12950
$ws.Range("F6").Value  = 7263
$ws.Range("F8").Value  = 1025
$ws.Range("F9").Value  = 160
$ws.Range("F10").Value = 373
$ws.Range("F11").Value = 1055
$ws.Range("F12").Value = 18
$ws.Range("F16").Value = 277
$ws.Range("F19").Value = 292
$ws.Range("F20").Value = 324
$ws.Range("F25").Value = 264
$ws.Range("F27").Value = 5322
$ws.Range("F32").Value = 331
$ws.Range("F34").Value = 2040
$ws.Range("F35").Value = 105
$ws.Range("F37").Value = 1404
$ws.Range("F43").Value = 82
$ws.Range("F47").Value = 3758
$ws.Range("F49").Value = 49
